$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 33
$ws_ALC.Range("H33").Value = 353.5
$ws_ALC.Range("I33").Value = 342.23077
$ws_ALC.Range("J33").Value = 500
$ws_ALC.Range("K33").Value = 342.23077
$ws_ALC.Range("L33").Value = 500
$ws_ALC.Range("M33").Value = -113.23077
$ws_ALC.Range("N33").Value = -958

# ALC row 45
$ws_ALC.Range("H45").Value = 17093.334
$ws_ALC.Range("J45").Value = 0
$ws_ALC.Range("L45").Value = 0
$ws_ALC.Range("N45").ClearContents()

# ALC row 100
$ws_ALC.Range("H100").Value = 5682.636
$ws_ALC.Range("I100").Value = 5752.25
$ws_ALC.Range("J100").Value = 5642.857
$ws_ALC.Range("K100").Value = 5752.25
$ws_ALC.Range("L100").Value = 5642.857
$ws_ALC.Range("M100").Value = -5211.25
$ws_ALC.Range("N100").Value = -6724.857

# BSM row 88
$ws_BSM.Range("H88").Value = 38000
$ws_BSM.Range("J88").Value = 38000
$ws_BSM.Range("L88").Value = 38000
$ws_BSM.Range("N88").Value = -38812

# BSM row 91
$ws_BSM.Range("H91").Value = 38000
$ws_BSM.Range("J91").Value = 38000
$ws_BSM.Range("L91").Value = 38000
$ws_BSM.Range("N91").Value = -40808

# BSM row 107
$ws_BSM.Range("H107").Value = 65887.125
$ws_BSM.Range("I107").Value = 74971
$ws_BSM.Range("J107").Value = 2300
$ws_BSM.Range("K107").Value = 74971
$ws_BSM.Range("L107").Value = 2300
$ws_BSM.Range("M107").Value = -73051
$ws_BSM.Range("N107").Value = -6140

# CRP row 31
$ws_CRP.Range("H31").Value = 2286.6445
$ws_CRP.Range("I31").Value = 2286.875
$ws_CRP.Range("J31").Value = 2286.3809
$ws_CRP.Range("K31").Value = 2286.875
$ws_CRP.Range("L31").Value = 2286.3809
$ws_CRP.Range("M31").Value = -1991.875
$ws_CRP.Range("N31").Value = -2876.3809

# CRP row 34
$ws_CRP.Range("H34").Value = 2286.6445
$ws_CRP.Range("I34").Value = 2286.875
$ws_CRP.Range("J34").Value = 2286.3809
$ws_CRP.Range("K34").Value = 2286.875
$ws_CRP.Range("L34").Value = 2286.3809
$ws_CRP.Range("M34").Value = -2084.875
$ws_CRP.Range("N34").Value = -2690.3809

# CRP row 107
$ws_CRP.Range("H107").Value = 380
$ws_CRP.Range("I107").Value = 370.48
$ws_CRP.Range("J107").Value = 399.83334
$ws_CRP.Range("K107").Value = 370.48
$ws_CRP.Range("L107").Value = 399.83334
$ws_CRP.Range("M107").Value = 1549.52
$ws_CRP.Range("N107").Value = -4239.83334

# CUL row 17
$ws_CUL.Range("H17").Value = 6000

# CUL row 34
$ws_CUL.Range("H34").Value = 995
$ws_CUL.Range("I34").Value = 392
$ws_CUL.Range("J34").Value = 2000
$ws_CUL.Range("K34").Value = 1176
$ws_CUL.Range("L34").Value = 6000
$ws_CUL.Range("M34").Value = -1092
$ws_CUL.Range("N34").Value = -6168

# CUL row 39
$ws_CUL.Range("H39").Value = 5695.2
$ws_CUL.Range("I39").Value = 0
$ws_CUL.Range("J39").Value = 5695.2
$ws_CUL.Range("K39").Value = 0
$ws_CUL.Range("L39").Value = 17085.6
$ws_CUL.Range("N39").Value = -17673.6
$ws_CUL.Range("M39").ClearContents()

# CUL row 55
$ws_CUL.Range("H55").Value = 3758.8462
$ws_CUL.Range("J55").Value = 3863.9167
$ws_CUL.Range("L55").Value = 11591.7501
$ws_CUL.Range("N55").Value = -11945.7501

# CUL row 68
$ws_CUL.Range("H68").Value = 752.0925999999999
$ws_CUL.Range("I68").Value = 565.94116
$ws_CUL.Range("J68").Value = 1068.55
$ws_CUL.Range("K68").Value = 1697.82348
$ws_CUL.Range("L68").Value = 3205.65
$ws_CUL.Range("M68").Value = -886.82348
$ws_CUL.Range("N68").Value = -4827.65

# CUL row 71
$ws_CUL.Range("H71").Value = 752.0925999999999
$ws_CUL.Range("I71").Value = 565.94116
$ws_CUL.Range("J71").Value = 1068.55
$ws_CUL.Range("K71").Value = 5093.47044
$ws_CUL.Range("L71").Value = 9616.949999999999
$ws_CUL.Range("M71").Value = -1037.47044
$ws_CUL.Range("N71").Value = -17728.95

# CUL row 82
$ws_CUL.Range("H82").Value = 24004
$ws_CUL.Range("J82").Value = 24004
$ws_CUL.Range("L82").Value = 72012
$ws_CUL.Range("N82").Value = -72824

# CUL row 85
$ws_CUL.Range("H85").Value = 24004
$ws_CUL.Range("J85").Value = 24004
$ws_CUL.Range("L85").Value = 72012
$ws_CUL.Range("N85").Value = -74820

# CUL row 107
$ws_CUL.Range("H107").Value = 773.2771
$ws_CUL.Range("I107").Value = 764.8108
$ws_CUL.Range("J107").Value = 780.087
$ws_CUL.Range("K107").Value = 2294.4324
$ws_CUL.Range("L107").Value = 2340.261
$ws_CUL.Range("M107").Value = -374.4323999999997
$ws_CUL.Range("N107").Value = -6180.261

# GSM row 107
$ws_GSM.Range("H107").Value = 1383.6666
$ws_GSM.Range("I107").Value = 1575.5
$ws_GSM.Range("J107").Value = 1000
$ws_GSM.Range("K107").Value = 1575.5
$ws_GSM.Range("L107").Value = 1000
$ws_GSM.Range("M107").Value = 344.5
$ws_GSM.Range("N107").Value = -4840

# LTW row 61
$ws_LTW.Range("H61").Value = 2853.3845
$ws_LTW.Range("I61").Value = 2601.2727
$ws_LTW.Range("K61").Value = 2601.2727
$ws_LTW.Range("M61").Value = -2399.2727

# LTW row 108
$ws_LTW.Range("H108").Value = 37000
$ws_LTW.Range("J108").Value = 37000
$ws_LTW.Range("L108").Value = 37000
$ws_LTW.Range("N108").Value = -44680

# LTW row 113
$ws_LTW.Range("H113").Value = 2853.3845
$ws_LTW.Range("I113").Value = 2601.2727
$ws_LTW.Range("K113").Value = 2601.2727
$ws_LTW.Range("M113").Value = -431.2727

# LTW row 122
$ws_LTW.Range("H122").Value = 90003450
$ws_LTW.Range("I122").Value = 83336420
$ws_LTW.Range("K122").Value = 250009260
$ws_LTW.Range("M122").Value = -250006810

# LTW row 129
$ws_LTW.Range("H129").Value = 40429
$ws_LTW.Range("J129").Value = 40429
$ws_LTW.Range("L129").Value = 40429
$ws_LTW.Range("N129").Value = -50429

# LTW row 136
$ws_LTW.Range("H136").Value = 1291.1333
$ws_LTW.Range("I136").Value = 863.9167
$ws_LTW.Range("K136").Value = 2591.7501
$ws_LTW.Range("M136").Value = -41.7501000000002

# WVR row 13
$ws_WVR.Range("H13").Value = 1916.6666
$ws_WVR.Range("J13").Value = 300
$ws_WVR.Range("L13").Value = 300
$ws_WVR.Range("N13").Value = -580

# WVR row 70
$ws_WVR.Range("H70").Value = 152475
$ws_WVR.Range("I70").Value = 136633.33
$ws_WVR.Range("K70").Value = 136633.33
$ws_WVR.Range("M70").Value = -136318.33

# WVR row 73
$ws_WVR.Range("H73").Value = 152475
$ws_WVR.Range("I73").Value = 136633.33
$ws_WVR.Range("K73").Value = 136633.33
$ws_WVR.Range("M73").Value = -135541.33

# WVR row 107
$ws_WVR.Range("H107").Value = 428.4
$ws_WVR.Range("I107").Value = 400.5
$ws_WVR.Range("J107").Value = 447
$ws_WVR.Range("K107").Value = 1201.5
$ws_WVR.Range("L107").Value = 1341
$ws_WVR.Range("M107").Value = 718.5
$ws_WVR.Range("N107").Value = -5181

# WVR row 129
$ws_WVR.Range("H129").Value = 0
$ws_WVR.Range("J129").Value = 0
$ws_WVR.Range("L129").Value = 0
$ws_WVR.Range("N129").ClearContents()

# WVR row 136
$ws_WVR.Range("H136").Value = 2227.4333
$ws_WVR.Range("I136").Value = 2088
$ws_WVR.Range("J136").Value = 2552.7778
$ws_WVR.Range("K136").Value = 6264
$ws_WVR.Range("L136").Value = 7658.3334
$ws_WVR.Range("M136").Value = -3714
$ws_WVR.Range("N136").Value = -12758.3334
